$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 20 (pushes old rows 20-26 down to 23-29).
$ws.Rows.Item(20).Resize(3).Insert()

# --- Row 19: "physical_interfaces" section header row now starts the new
#     merged block B19:B22 and its data columns describe "Sleep Unit".
$ws.Cells.Item(19, 3).Value = "Sleep Unit"
$ws.Cells.Item(19, 4).Value = "Needs capture"
$ws.Cells.Item(19, 5).Value = ""
$ws.Cells.Item(19, 6).Value = ""

# --- Row 20 (new): "Configuration and Control"
$ws.Cells.Item(20, 3).Value = "Configuration and Control"
$ws.Cells.Item(20, 4).Value = "Needs capture"
$ws.Cells.Item(20, 5).Value = ""
$ws.Cells.Item(20, 6).Value = ""
$ws.Cells.Item(20, 3).WrapText = $true
$ws.Rows.Item(20).RowHeight = 23.45

# --- Row 21 (new): "APU"
$ws.Cells.Item(21, 3).Value = "APU"
$ws.Cells.Item(21, 4).Value = "Out of scope"
$ws.Cells.Item(21, 5).Value = ""
$ws.Cells.Item(21, 6).Value = ""
$ws.Cells.Item(21, 3).WrapText = $true

# --- Row 22 (new): the original "OBI" row content that used to live on row 19.
$ws.Cells.Item(22, 3).Value = "OBI"
$ws.Cells.Item(22, 4).Value = "Ready for Review"
$ws.Cells.Item(22, 5).Value = "Arjan B"
$ws.Cells.Item(22, 6).Value = ""

# Copy the row-19 style down into the 3 new rows so B/C/D/E/F formatting
# (borders, fonts, alignment) matches the rest of the merged section.
$ws.Range("B19:F19").Copy()
$ws.Range("B20:F22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-apply the text/wrap specific formatting that PasteSpecial(formats) may
# have clobbered back to their intended values.
$ws.Cells.Item(20, 3).WrapText = $true
$ws.Rows.Item(20).RowHeight = 23.45
$ws.Cells.Item(21, 3).WrapText = $true

# B19 becomes part of the merged block, matching the style used by the
# other merged section headers (B13, B16, B20-before-edit, ...).
$ws.Cells.Item(19, 2).Style = $ws.Cells.Item(13, 2).Style

# Re-create the merges: B19:B22 (physical_interfaces) and B23:B29 /
# F23:F29 (xpulp_instruction_extensions), replacing the old B20:B26/F20:F26
# pair that the row-insert above already shifted down to B23:B29/F23:F29.
$ws.Range("B19:B22").Merge()

Write-Output "done"
